# Generate Report for Handoff
# Promotes the four "Ready for handoff" entries (180a04e6.., 6f39b068..,
# 88f8d741.., cdc67cd2..) to the just-handed-off state: Priority moves from
# "low" to "ht" on both locale sheets, and the handoff timestamps are
# refreshed.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# zh-cn: Priority low -> ht, and a fresh "Latest Handoff Datetime"
foreach ($r in 4..7) {
    $zhcn.Cells.Item($r, 5).Value = "ht"
    $zhcn.Cells.Item($r, 8).Value = "2016-09-03 22:34:25"
}

# de-de: Priority low -> ht, and a fresh "Latest Handoff Datetime" that
# matches the refreshed Overview generate-date below
foreach ($r in 4..7) {
    $dede.Cells.Item($r, 5).Value = "ht"
    $dede.Cells.Item($r, 8).Value = "2016-09-03 22:34:29"
}

# Overview: refresh "Latest HO Xliff Generate Date" for the same four rows
foreach ($r in 4..7) {
    $overview.Cells.Item($r, 7).Value = "2016-09-03 22:34:29"
}
